# Ccl25-Ackr4: refresh sheet with updated TPM-based NATMI output.
# New ligand/receptor expression values for existing Sending x Target cluster
# pairs, plus four new rows for the "Resolving-Mac" sending cluster (rows 14-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @('ECs', 'Ccl25', 'Ackr4', 'ECs', 3, 1, 4.055275000000001, 12.165825, 0.1947228515851206, 0.1947228515851206, 1, 0.3333333333333333, 0.1143813333333333, 0.343144, 0.128300337591142, 0.1283003375911419, 0.4638477615333335, 4.174629853800001, 0.02498300759508081, 0.0249830075950808),
    @('ECs', 'Ccl25', 'Ackr4', 'FAPs', 3, 1, 4.055275000000001, 12.165825, 0.1947228515851206, 0.1947228515851206, 3, 1, 0.7200953333333334, 2.160286, 0.8077233554817153, 0.8077233554817151, 2.920184602883334, 26.28166142595001, 0.1572821950713017, 0.1572821950713016),
    @('ECs', 'Ccl25', 'Ackr4', 'MuSCs', 3, 1, 4.055275000000001, 12.165825, 0.1947228515851206, 0.1947228515851206, 1, 0.3333333333333333, 0.03357866666666667, 0.100736, 0.03766483694187069, 0.03766483694187069, 0.1361707274666667, 1.2255365472, 0.007334204453809654, 0.007334204453809653),
    @('ECs', 'Ccl25', 'Ackr4', 'Resolving-Mac', 3, 1, 4.055275000000001, 12.165825, 0.1947228515851206, 0.1947228515851206, 1, 0.3333333333333333, 0.023457, 0.070371, 0.02631146998527222, 0.02631146998527222, 0.09512458567500003, 0.8561212710750001, 0.005123444464928518, 0.005123444464928516),
    @('FAPs', 'Ccl25', 'Ackr4', 'ECs', 3, 1, 6.542812333333333, 19.628437, 0.3141673684110111, 0.3141673684110111, 1, 0.3333333333333333, 0.1143813333333333, 0.343144, 0.128300337591142, 0.1283003375911419, 0.7483755984364444, 6.735380385927999, 0.04030777942725339, 0.04030777942725339),
    @('FAPs', 'Ccl25', 'Ackr4', 'FAPs', 3, 1, 6.542812333333333, 19.628437, 0.3141673684110111, 0.3141673684110111, 3, 1, 0.7200953333333334, 2.160286, 0.8077233554817153, 0.8077233554817151, 4.711448628109111, 42.403037652982, 0.2537603209958021, 0.2537603209958021),
    @('FAPs', 'Ccl25', 'Ackr4', 'MuSCs', 3, 1, 6.542812333333333, 19.628437, 0.3141673684110111, 0.3141673684110111, 1, 0.3333333333333333, 0.03357866666666667, 0.100736, 0.03766483694187069, 0.03766483694187069, 0.2196989144035555, 1.977290229632, 0.01183306270365735, 0.01183306270365735),
    @('FAPs', 'Ccl25', 'Ackr4', 'Resolving-Mac', 3, 1, 6.542812333333333, 19.628437, 0.3141673684110111, 0.3141673684110111, 1, 0.3333333333333333, 0.023457, 0.070371, 0.02631146998527222, 0.02631146998527222, 0.153474748903, 1.381272740127, 0.008266205284298278, 0.008266205284298278),
    @('MuSCs', 'Ccl25', 'Ackr4', 'ECs', 3, 1, 6.7998, 20.3994, 0.3265071903159472, 0.3265071903159472, 1, 0.3333333333333333, 0.1143813333333333, 0.343144, 0.128300337591142, 0.1283003375911419, 0.7777701904000001, 6.9999317136, 0.04189098274347126, 0.04189098274347126),
    @('MuSCs', 'Ccl25', 'Ackr4', 'FAPs', 3, 1, 6.7998, 20.3994, 0.3265071903159472, 0.3265071903159472, 3, 1, 0.7200953333333334, 2.160286, 0.8077233554817153, 0.8077233554817151, 4.896504247600001, 44.0685382284, 0.2637274833509039, 0.2637274833509038),
    @('MuSCs', 'Ccl25', 'Ackr4', 'MuSCs', 3, 1, 6.7998, 20.3994, 0.3265071903159472, 0.3265071903159472, 1, 0.3333333333333333, 0.03357866666666667, 0.100736, 0.03766483694187069, 0.03766483694187069, 0.2283282176, 2.0549539584, 0.01229784008359849, 0.01229784008359849),
    @('MuSCs', 'Ccl25', 'Ackr4', 'Resolving-Mac', 3, 1, 6.7998, 20.3994, 0.3265071903159472, 0.3265071903159472, 1, 0.3333333333333333, 0.023457, 0.070371, 0.02631146998527222, 0.02631146998527222, 0.1595029086, 1.4355261774, 0.00859088413797361, 0.008590884137973608),
    @('Resolving-Mac', 'Ccl25', 'Ackr4', 'ECs', 3, 1, 3.427994, 10.283982, 0.164602589687921, 0.164602589687921, 1, 0.3333333333333333, 0.1143813333333333, 0.343144, 0.128300337591142, 0.1283003375911419, 0.3920985243786667, 3.528886719408, 0.02111856782533648, 0.02111856782533648),
    @('Resolving-Mac', 'Ccl25', 'Ackr4', 'FAPs', 3, 1, 3.427994, 10.283982, 0.164602589687921, 0.164602589687921, 3, 1, 0.7200953333333334, 2.160286, 0.8077233554817153, 0.8077233554817151, 2.468482482094667, 22.216342338852, 0.1329533560637075, 0.1329533560637075),
    @('Resolving-Mac', 'Ccl25', 'Ackr4', 'MuSCs', 3, 1, 3.427994, 10.283982, 0.164602589687921, 0.164602589687921, 1, 0.3333333333333333, 0.03357866666666667, 0.100736, 0.03766483694187069, 0.03766483694187069, 0.1151074678613333, 1.035967210752, 0.006199729700805189, 0.006199729700805189),
    @('Resolving-Mac', 'Ccl25', 'Ackr4', 'Resolving-Mac', 3, 1, 3.427994, 10.283982, 0.164602589687921, 0.164602589687921, 1, 0.3333333333333333, 0.023457, 0.070371, 0.02631146998527222, 0.02631146998527222, 0.080410455258, 0.7236940973220001, 0.004330936098071812, 0.004330936098071811)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowVals[$j]
    }
}

Write-Output "Done: wrote $($data.Count) rows starting at row $startRow"
